# Review and define analysis section fields
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) AnalysisFields sheet: populate with a new summary table (Table4) that
#    reviews the "average_axis" / "average_diameter" / "average_weight"
#    fields, mirroring the layout already used on ExamplesFields.
# ---------------------------------------------------------------------------
$afs = $wb.Worksheets.Item("AnalysisFields")

$afs.Range("A1").Value = "field_name"
$afs.Range("B1").Value = "data_type"
$afs.Range("C1").Value = "example_field_value"
$afs.Range("D1").Value = "example_field_uri"

$afs.Range("A2").Value = "average_axis"
$afs.Range("B2").Value = "NUMERIC"
$afs.Range("C2").Value = 5
$afs.Range("D2").Value = "http://numismatics.org/ocre/id/ric.1(2).aug.1A"

$afs.Range("A3").Value = "average_diameter"
$afs.Range("B3").Value = "NUMERIC"
$afs.Range("C3").Value = 13.66
$afs.Range("D3").Value = "http://numismatics.org/ocre/id/ric.1(2).aug.1A"

$afs.Range("A4").Value = "average_weight"
$afs.Range("B4").Value = "NUMERIC"
$afs.Range("C4").Value = 1.61
$afs.Range("D4").Value = "http://numismatics.org/ocre/id/ric.1(2).aug.1A"

$analysisTable = $afs.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $afs.Range("A1:D4"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$analysisTable.Name = "Table4"
$analysisTable.TableStyle = "TableStyleMedium2"
$analysisTable.ShowTotals = $true

$afs.Range("A5").Value = "Total"
$afs.Range("D5").Formula = "=SUBTOTAL(103,Table4[example_field_uri])"
$analysisTable.ListColumns.Item("example_field_uri").TotalsCalculation = [Microsoft.Office.Interop.Excel.XlTotalsCalculation]::xlTotalsCalculationCount

$afs.Columns.Item(1).ColumnWidth = 16.5
$afs.Columns.Item(2).ColumnWidth = 16.5
$afs.Columns.Item(3).ColumnWidth = 20.6640625
$afs.Columns.Item(4).ColumnWidth = 40.6640625

$afs.Range("C6").Select()

# ---------------------------------------------------------------------------
# 2) ExamplesFields sheet: turn on the totals row for the existing table so
#    it reports a count of the first_example_uri examples reviewed.
# ---------------------------------------------------------------------------
$efs = $wb.Worksheets.Item("ExamplesFields")
$examplesTable = $efs.ListObjects.Item(1)
$examplesTable.ShowTotals = $true

$efs.Range("A9").Value = "Total"
$efs.Range("E9").Formula = "=SUBTOTAL(103," + $examplesTable.Name + "[first_example_uri])"
$examplesTable.ListColumns.Item("first_example_uri").TotalsCalculation = [Microsoft.Office.Interop.Excel.XlTotalsCalculation]::xlTotalsCalculationCount

# ---------------------------------------------------------------------------
# 3) Collections sheet: selection moved after reviewing the table contents.
# ---------------------------------------------------------------------------
$cols = $wb.Worksheets.Item("Collections")
$cols.Range("D57").Select()

# ---------------------------------------------------------------------------
# 4) CollectionsAndIIIF sheet: header row reverts to the default (unbolded)
#    style now that the custom header style is no longer used elsewhere.
# ---------------------------------------------------------------------------
$ciiif = $wb.Worksheets.Item("CollectionsAndIIIF")
$ciiif.Range("A1:E1").Style = "Normal"

# ---------------------------------------------------------------------------
# Make AnalysisFields the active / visible tab, matching the reviewed state.
# ---------------------------------------------------------------------------
$afs.Activate()
